$wb = $excel.ActiveWorkbook

# --- AgentSettings: "Mixed" -> "Averaged" (C15) ---
$ws1 = $wb.Worksheets.Item("AgentSettings")
$ws1.Range("C15").Value = "Averaged"

# --- StrategySpaceInReference5: StrategicFactors list updated (D2) ---
$ws2 = $wb.Worksheets.Item("StrategySpaceInReference5")
$ws2.Range("D2").Value = "StrategicFactors: 1,1.5,2,2.5,3,3.5,4,4.5,5,5.5,6,6.5,7,7.5,8,8.5,9,9.5,10"

# --- StrategySpaceInReference7: StrategicFactors / CostShape lists updated (D2, D3) ---
$ws3 = $wb.Worksheets.Item("StrategySpaceInReference7")
$ws3.Range("D2").Value = "StrategicFactors : 1,1.25,1.5,1.75,2,2.25,2.5,2.75,3,3.25,3.5,3.75,4,4.25,4.5,4.75,5,5.25,5.5"
$ws3.Range("D3").Value = "CostShape: [1,1.1,1.2,1.3,1.4,1.5,1.6,1.7,1.8,1.9]"

# --- Selection / active sheet bookkeeping to mirror the recorded view state ---
$ws2.Range("D21").Select()
$ws3.Range("D7").Select()
$ws1.Select()
$ws1.Range("C25").Select()
